# Atualização de dados 22/01/2024: SKU e CNPJ
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New SKU rows appended below the existing data (A20:A27). Copy the
# formatting of the last existing data row (A19) down onto the new rows so
# they pick up the same style (centered integer, thin border) already used
# by the rest of the column, then fill in the new values.
$newValues = @(5255, 5279, 5262, 3342, 341, 5354, 5361, 5378)

$startRow = 20
$endRow = $startRow + $newValues.Length - 1

$ws.Range("A19").Copy()
$ws.Range("A" + $startRow + ":A" + $endRow).PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# Update the view: scrolled so row 2 is the top-left visible row, with the
# active selection on F15.
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("F15").Select()
